$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2..17 (columns A..T), reflecting the
# re-run of the NATMI analysis after adding the "Resolving-Mac" cluster
# to the dataset (new sending/target cluster + refreshed TPM-derived
# statistics for all existing sending/target cluster combinations).
$data = @(
    ,@("ECs","Vtn","Itgav","ECs",3,1,2.767552,8.302655999999999,0.04706493447833917,0.04706493447833917,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,26.08100137282133,234.729012355392,0.003131859014957021,0.003131859014957021)
    ,@("ECs","Vtn","Itgav","FAPs",3,1,2.767552,8.302655999999999,0.04706493447833917,0.04706493447833917,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,140.0175535283627,1260.157981755264,0.01681358898002286,0.01681358898002286)
    ,@("ECs","Vtn","Itgav","MuSCs",3,1,2.767552,8.302655999999999,0.04706493447833917,0.04706493447833917,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,74.28320639965867,668.548857596928,0.008920076583605406,0.008920076583605408)
    ,@("ECs","Vtn","Itgav","Resolving-Mac",3,1,2.767552,8.302655999999999,0.04706493447833917,0.04706493447833917,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,151.558174334528,1364.023569010752,0.01819940989975388,0.01819940989975388)
    ,@("FAPs","Vtn","Itgav","ECs",3,1,21.05317333333333,63.15952,0.3580298485789791,0.3580298485789791,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,198.4019966414045,1785.61796977264,0.02382451014378511,0.02382451014378511)
    ,@("FAPs","Vtn","Itgav","FAPs",3,1,21.05317333333333,63.15952,0.3580298485789791,0.3580298485789791,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,1065.134033305209,9586.20629974688,0.1279034334862884,0.1279034334862884)
    ,@("FAPs","Vtn","Itgav","MuSCs",3,1,21.05317333333333,63.15952,0.3580298485789791,0.3580298485789791,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,565.0832288201956,5085.74905938176,0.06785632879210668,0.06785632879210669)
    ,@("FAPs","Vtn","Itgav","Resolving-Mac",3,1,21.05317333333333,63.15952,0.3580298485789791,0.3580298485789791,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,1152.925225740427,10376.32703166384,0.1384455761567989,0.1384455761567989)
    ,@("MuSCs","Vtn","Itgav","ECs",3,1,34.97741266666667,104.932238,0.5948251867999219,0.5948251867999219,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,329.6219719727296,2966.597747754567,0.03958166826855355,0.03958166826855355)
    ,@("MuSCs","Vtn","Itgav","FAPs",3,1,34.97741266666667,104.932238,0.5948251867999219,0.5948251867999219,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,1769.59701221102,15926.37310989917,0.2124967625403167,0.2124967625403167)
    ,@("MuSCs","Vtn","Itgav","MuSCs",3,1,34.97741266666667,104.932238,0.5948251867999219,0.5948251867999219,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,938.8204320800606,8449.383888720546,0.1127354426160869,0.1127354426160869)
    ,@("MuSCs","Vtn","Itgav","Resolving-Mac",3,1,34.97741266666667,104.932238,0.5948251867999219,0.5948251867999219,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,1915.451925277428,17239.06732749685,0.2300113133749647,0.2300113133749647)
    ,@("Resolving-Mac","Vtn","Itgav","ECs",1,0.3333333333333333,0.004706,0.014118,0.00008003014275976175,0.00008003014275976175,3,1,9.423852333333334,28.271557,0.06654336290212845,0.06654336290212845,0.04434864908066667,0.399137841726,0.000005325474832771974,0.000005325474832771974)
    ,@("Resolving-Mac","Vtn","Itgav","FAPs",1,0.3333333333333333,0.004706,0.014118,0.00008003014275976175,0.00008003014275976175,3,1,50.59256466666667,151.777694,0.3572423751649123,0.3572423751649123,0.2380886093213334,2.142797483892,0.0000285901582842843,0.0000285901582842843)
    ,@("Resolving-Mac","Vtn","Itgav","MuSCs",1,0.3333333333333333,0.004706,0.014118,0.00008003014275976175,0.00008003014275976175,3,1,26.84076266666667,80.522288,0.1895270158659356,0.1895270158659356,0.1263126291093333,1.136813661984,0.00001516787413658245,0.00001516787413658246)
    ,@("Resolving-Mac","Vtn","Itgav","Resolving-Mac",1,0.3333333333333333,0.004706,0.014118,0.00008003014275976175,0.00008003014275976175,3,1,54.762539,164.287617,0.3866872460670236,0.3866872460670236,0.257712508534,2.319412576806,0.00003094663550612302,0.00003094663550612302)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
